$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H64").Value = 4358.6113
$ws_ALC.Range("J64").Value = 4430.3335
$ws_ALC.Range("L64").Value = 4430.3335
$ws_ALC.Range("N64").Value = -4926.3335
$ws_ALC.Range("H67").Value = 4358.6113
$ws_ALC.Range("J67").Value = 4430.3335
$ws_ALC.Range("L67").Value = 4430.3335
$ws_ALC.Range("N67").Value = -6146.3335
$ws_ALC.Range("H74").Value = 3933.5
$ws_ALC.Range("I74").Value = 3867
$ws_ALC.Range("K74").Value = 3867
$ws_ALC.Range("M74").Value = -2931
$ws_ALC.Range("H76").Value = 3540.5
$ws_ALC.Range("I76").Value = 3411.111
$ws_ALC.Range("J76").Value = 3646.3635
$ws_ALC.Range("K76").Value = 3411.111
$ws_ALC.Range("L76").Value = 3646.3635
$ws_ALC.Range("M76").Value = -3096.111
$ws_ALC.Range("N76").Value = -4276.363499999999
$ws_ALC.Range("H77").Value = 3933.5
$ws_ALC.Range("I77").Value = 3867
$ws_ALC.Range("K77").Value = 19335
$ws_ALC.Range("M77").Value = -14655
$ws_ALC.Range("H79").Value = 3540.5
$ws_ALC.Range("I79").Value = 3411.111
$ws_ALC.Range("J79").Value = 3646.3635
$ws_ALC.Range("K79").Value = 3411.111
$ws_ALC.Range("L79").Value = 3646.3635
$ws_ALC.Range("M79").Value = -2319.111
$ws_ALC.Range("N79").Value = -5830.363499999999
$ws_ALC.Range("H82").Value = 0
$ws_ALC.Range("I82").Value = 0
$ws_ALC.Range("K82").Value = 0
$ws_ALC.Range("M82").ClearContents()
$ws_ALC.Range("H85").Value = 0
$ws_ALC.Range("I85").Value = 0
$ws_ALC.Range("K85").Value = 0
$ws_ALC.Range("M85").ClearContents()
$ws_ALC.Range("H86").Value = 2751.7693
$ws_ALC.Range("I86").Value = 1503
$ws_ALC.Range("J86").Value = 2855.8333
$ws_ALC.Range("K86").Value = 1503
$ws_ALC.Range("L86").Value = 2855.8333
$ws_ALC.Range("M86").Value = -380
$ws_ALC.Range("N86").Value = -5101.8333
$ws_ALC.Range("H87").Value = 30331.666
$ws_ALC.Range("J87").Value = 30331.666
$ws_ALC.Range("L87").Value = 30331.666
$ws_ALC.Range("N87").Value = -32827.666
$ws_ALC.Range("H88").Value = 23523.555
$ws_ALC.Range("I88").Value = 976
$ws_ALC.Range("J88").Value = 29965.715
$ws_ALC.Range("K88").Value = 976
$ws_ALC.Range("L88").Value = 29965.715
$ws_ALC.Range("M88").Value = -570
$ws_ALC.Range("N88").Value = -30777.715
$ws_ALC.Range("H89").Value = 2751.7693
$ws_ALC.Range("I89").Value = 1503
$ws_ALC.Range("J89").Value = 2855.8333
$ws_ALC.Range("K89").Value = 7515
$ws_ALC.Range("L89").Value = 14279.1665
$ws_ALC.Range("M89").Value = -1899
$ws_ALC.Range("N89").Value = -25511.1665
$ws_ALC.Range("H90").Value = 30331.666
$ws_ALC.Range("J90").Value = 30331.666
$ws_ALC.Range("L90").Value = 90994.99800000001
$ws_ALC.Range("N90").Value = -103474.998
$ws_ALC.Range("H91").Value = 23523.555
$ws_ALC.Range("I91").Value = 976
$ws_ALC.Range("J91").Value = 29965.715
$ws_ALC.Range("K91").Value = 976
$ws_ALC.Range("L91").Value = 29965.715
$ws_ALC.Range("M91").Value = 428
$ws_ALC.Range("N91").Value = -32773.715
$ws_ALC.Range("H94").Value = 3000
$ws_ALC.Range("I94").Value = 3000
$ws_ALC.Range("K94").Value = 3000
$ws_ALC.Range("M94").Value = -2549
$ws_ALC.Range("H131").Value = 877.5
$ws_ALC.Range("I131").Value = 760
$ws_ALC.Range("J131").Value = 995
$ws_ALC.Range("K131").Value = 2280
$ws_ALC.Range("L131").Value = 2985
$ws_ALC.Range("M131").Value = 2760
$ws_ALC.Range("N131").Value = -13065
$ws_ALC.Range("H141").Value = 2493.125
$ws_ALC.Range("I141").Value = 2114.6155
$ws_ALC.Range("J141").Value = 4133.3335
$ws_ALC.Range("K141").Value = 6343.8465
$ws_ALC.Range("L141").Value = 12400.0005
$ws_ALC.Range("M141").Value = -1163.8465
$ws_ALC.Range("N141").Value = -22760.0005

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 2799.9
$ws_ARM.Range("I45").Value = 3000
$ws_ARM.Range("J45").Value = 2333
$ws_ARM.Range("K45").Value = 3000
$ws_ARM.Range("L45").Value = 2333
$ws_ARM.Range("M45").Value = -2623
$ws_ARM.Range("N45").Value = -3087
$ws_ARM.Range("H74").Value = 7201204.5
$ws_ARM.Range("I74").Value = 10914166
$ws_ARM.Range("K74").Value = 10914166
$ws_ARM.Range("M74").Value = -10913292
$ws_ARM.Range("H77").Value = 7201204.5
$ws_ARM.Range("I77").Value = 10914166
$ws_ARM.Range("K77").Value = 54570830
$ws_ARM.Range("M77").Value = -54566462
$ws_ARM.Range("H122").Value = 2060.6
$ws_ARM.Range("I122").Value = 2101
$ws_ARM.Range("K122").Value = 6303
$ws_ARM.Range("M122").Value = -3853
$ws_ARM.Range("H123").Value = 60000
$ws_ARM.Range("J123").Value = 60000
$ws_ARM.Range("L123").Value = 60000
$ws_ARM.Range("N123").Value = -69800

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H68").Value = 36900
$ws_BSM.Range("J68").Value = 36900
$ws_BSM.Range("L68").Value = 36900
$ws_BSM.Range("N68").Value = -38522
$ws_BSM.Range("H69").Value = 40295
$ws_BSM.Range("J69").Value = 40295
$ws_BSM.Range("L69").Value = 40295
$ws_BSM.Range("N69").Value = -41917
$ws_BSM.Range("H71").Value = 36900
$ws_BSM.Range("J71").Value = 36900
$ws_BSM.Range("L71").Value = 110700
$ws_BSM.Range("N71").Value = -118812
$ws_BSM.Range("H72").Value = 40295
$ws_BSM.Range("J72").Value = 40295
$ws_BSM.Range("L72").Value = 120885
$ws_BSM.Range("N72").Value = -128997
$ws_BSM.Range("H75").Value = 5071.3335
$ws_BSM.Range("I75").Value = 3607
$ws_BSM.Range("K75").Value = 3607
$ws_BSM.Range("M75").Value = -2671
$ws_BSM.Range("H76").Value = 36862.184
$ws_BSM.Range("J76").Value = 36862.184
$ws_BSM.Range("L76").Value = 36862.184
$ws_BSM.Range("N76").Value = -37492.184
$ws_BSM.Range("H78").Value = 5071.3335
$ws_BSM.Range("I78").Value = 3607
$ws_BSM.Range("K78").Value = 10821
$ws_BSM.Range("M78").Value = -6141
$ws_BSM.Range("H79").Value = 36862.184
$ws_BSM.Range("J79").Value = 36862.184
$ws_BSM.Range("L79").Value = 36862.184
$ws_BSM.Range("N79").Value = -39046.184
$ws_BSM.Range("H82").Value = 17887.166
$ws_BSM.Range("I82").Value = 6689.25
$ws_BSM.Range("J82").Value = 40283
$ws_BSM.Range("K82").Value = 6689.25
$ws_BSM.Range("L82").Value = 40283
$ws_BSM.Range("M82").Value = -6306.25
$ws_BSM.Range("N82").Value = -41049
$ws_BSM.Range("H85").Value = 17887.166
$ws_BSM.Range("I85").Value = 6689.25
$ws_BSM.Range("J85").Value = 40283
$ws_BSM.Range("K85").Value = 6689.25
$ws_BSM.Range("L85").Value = 40283
$ws_BSM.Range("M85").Value = -5363.25
$ws_BSM.Range("N85").Value = -42935
$ws_BSM.Range("H105").Value = 2495.9
$ws_BSM.Range("I105").Value = 2313.8
$ws_BSM.Range("J105").Value = 2678
$ws_BSM.Range("K105").Value = 2313.8
$ws_BSM.Range("L105").Value = 2678
$ws_BSM.Range("M105").Value = -566.8000000000002
$ws_BSM.Range("N105").Value = -6172

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H133").Value = 40748.383
$ws_CRP.Range("J133").Value = 40748.383
$ws_CRP.Range("L133").Value = 40748.383
$ws_CRP.Range("N133").Value = -45808.383
$ws_CRP.Range("H135").Value = 52000
$ws_CRP.Range("J135").Value = 52000
$ws_CRP.Range("L135").Value = 52000
$ws_CRP.Range("N135").Value = -62140
$ws_CRP.Range("H137").Value = 0
$ws_CRP.Range("J137").Value = 0
$ws_CRP.Range("L137").Value = 0
$ws_CRP.Range("N137").ClearContents()

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 26.533333
$ws_CUL.Range("I2").Value = 27.3
$ws_CUL.Range("J2").Value = 25
$ws_CUL.Range("K2").Value = 163.8
$ws_CUL.Range("L2").Value = 150
$ws_CUL.Range("M2").Value = -50.80000000000001
$ws_CUL.Range("N2").Value = -376
$ws_CUL.Range("H75").Value = 4000
$ws_CUL.Range("I75").Value = 0
$ws_CUL.Range("K75").Value = 0
$ws_CUL.Range("M75").ClearContents()
$ws_CUL.Range("H78").Value = 4000
$ws_CUL.Range("I78").Value = 0
$ws_CUL.Range("K78").Value = 0
$ws_CUL.Range("M78").ClearContents()
$ws_CUL.Range("H88").Value = 3362.5
$ws_CUL.Range("J88").Value = 3362.5
$ws_CUL.Range("L88").Value = 10087.5
$ws_CUL.Range("N88").Value = -10943.5
$ws_CUL.Range("H91").Value = 3362.5
$ws_CUL.Range("J91").Value = 3362.5
$ws_CUL.Range("L91").Value = 10087.5
$ws_CUL.Range("N91").Value = -13051.5
$ws_CUL.Range("H103").Value = 2419.75
$ws_CUL.Range("I103").Value = 422.77777
$ws_CUL.Range("J103").Value = 3617.9333
$ws_CUL.Range("K103").Value = 1268.33331
$ws_CUL.Range("L103").Value = 10853.7999
$ws_CUL.Range("M103").Value = -389.33331
$ws_CUL.Range("N103").Value = -12611.7999
$ws_CUL.Range("H107").Value = 986403.9399999999
$ws_CUL.Range("I107").Value = 932
$ws_CUL.Range("J107").Value = 1588636.9
$ws_CUL.Range("K107").Value = 2796
$ws_CUL.Range("L107").Value = 4765910.699999999
$ws_CUL.Range("M107").Value = -876
$ws_CUL.Range("N107").Value = -4769750.699999999
$ws_CUL.Range("H140").Value = 2055.0505
$ws_CUL.Range("I140").Value = 1299.2727
$ws_CUL.Range("J140").Value = 2270.987
$ws_CUL.Range("K140").Value = 3897.8181
$ws_CUL.Range("L140").Value = 6812.961
$ws_CUL.Range("N140").Value = -17172.961
$ws_CUL.Range("M140").Value = 1282.1819

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H126").Value = 1702.2273
$ws_GSM.Range("J126").Value = 2081.4
$ws_GSM.Range("L126").Value = 6244.200000000001
$ws_GSM.Range("N126").Value = -11184.2
